{"js": "// Replace the date line and each of the 100 multiplication-answer cells\n// with their updated values, in document order. Each old value is unique\n// in the document, so a simple search-and-replace per pair is safe.\nconst pairs = [\n  [\"2023-06-20 Tuesday\", \"2023-06-21 Wednesday\"],\n  [\"61\u00d747=2867\", \"17\u00d743=731\"],\n  [\"65\u00d752=3380\", \"83\u00d786=7138\"],\n  [\"25\u00d746=1150\", \"43\u00d786=3698\"],\n  [\"63\u00d746=2898\", \"13\u00d792=1196\"],\n  [\"14\u00d793=1302\", \"65\u00d781=5265\"],\n  [\"32\u00d743=1376\", \"87\u00d733=2871\"],\n  [\"96\u00d712=1152\", \"40\u00d714=560\"],\n  [\"14\u00d782=1148\", \"33\u00d728=924\"],\n  [\"86\u00d727=2322\", \"17\u00d766=1122\"],\n  [\"50\u00d737=1850\", \"37\u00d784=3108\"],\n  [\"10\u00d779=790\", \"96\u00d798=9408\"],\n  [\"44\u00d727=1188\", \"97\u00d748=4656\"],\n  [\"71\u00d725=1775\", \"59\u00d720=1180\"],\n  [\"12\u00d719=228\", \"46\u00d741=1886\"],\n  [\"70\u00d790=6300\", \"68\u00d759=4012\"],\n  [\"32\u00d765=2080\", \"94\u00d772=6768\"],\n  [\"76\u00d772=5472\", \"32\u00d710=320\"],\n  [\"45\u00d758=2610\", \"95\u00d780=7600\"],\n  [\"23\u00d778=1794\", \"70\u00d787=6090\"],\n  [\"60\u00d770=4200\", \"78\u00d788=6864\"],\n  [\"51\u00d774=3774\", \"82\u00d760=4920\"],\n  [\"43\u00d749=2107\", \"16\u00d794=1504\"],\n  [\"23\u00d734=782\", \"44\u00d711=484\"],\n  [\"80\u00d758=4640\", \"100\u00d732=3200\"],\n  [\"72\u00d767=4824\", \"79\u00d753=4187\"],\n  [\"18\u00d787=1566\", \"13\u00d758=754\"],\n  [\"78\u00d769=5382\", \"11\u00d762=682\"],\n  [\"75\u00d744=3300\", \"73\u00d724=1752\"],\n  [\"84\u00d721=1764\", \"38\u00d7100=3800\"],\n  [\"17\u00d759=1003\", \"16\u00d780=1280\"],\n  [\"25\u00d737=925\", \"37\u00d711=407\"],\n  [\"22\u00d788=1936\", \"74\u00d728=2072\"],\n  [\"14\u00d739=546\", \"81\u00d727=2187\"],\n  [\"79\u00d719=1501\", \"51\u00d731=1581\"],\n  [\"64\u00d757=3648\", \"34\u00d758=1972\"],\n  [\"53\u00d798=5194\", \"64\u00d799=6336\"],\n  [\"78\u00d712=936\", \"65\u00d742=2730\"],\n  [\"18\u00d788=1584\", \"56\u00d750=2800\"],\n  [\"97\u00d785=8245\", \"80\u00d731=2480\"],\n  [\"86\u00d763=5418\", \"42\u00d780=3360\"],\n  [\"12\u00d711=132\", \"94\u00d777=7238\"],\n  [\"29\u00d798=2842\", \"16\u00d742=672\"],\n  [\"40\u00d767=2680\", \"20\u00d723=460\"],\n  [\"93\u00d775=6975\", \"25\u00d715=375\"],\n  [\"82\u00d742=3444\", \"42\u00d771=2982\"],\n  [\"26\u00d756=1456\", \"88\u00d783=7304\"],\n  [\"73\u00d735=2555\", \"59\u00d778=4602\"],\n  [\"46\u00d768=3128\", \"37\u00d724=888\"],\n  [\"94\u00d768=6392\", \"24\u00d733=792\"],\n  [\"93\u00d716=1488\", \"98\u00d729=2842\"],\n  [\"64\u00d714=896\", \"93\u00d765=6045\"],\n  [\"16\u00d747=752\", \"38\u00d758=2204\"],\n  [\"52\u00d792=4784\", \"89\u00d776=6764\"],\n  [\"83\u00d733=2739\", \"69\u00d747=3243\"],\n  [\"33\u00d766=2178\", \"36\u00d797=3492\"],\n  [\"15\u00d732=480\", \"30\u00d716=480\"],\n  [\"66\u00d789=5874\", \"39\u00d769=2691\"],\n  [\"19\u00d743=817\", \"48\u00d799=4752\"],\n  [\"24\u00d776=1824\", \"85\u00d728=2380\"],\n  [\"20\u00d777=1540\", \"10\u00d735=350\"],\n  [\"28\u00d721=588\", \"65\u00d755=3575\"],\n  [\"27\u00d719=513\", \"56\u00d772=4032\"],\n  [\"20\u00d788=1760\", \"20\u00d718=360\"],\n  [\"16\u00d734=544\", \"10\u00d798=980\"],\n  [\"87\u00d732=2784\", \"58\u00d749=2842\"],\n  [\"40\u00d749=1960\", \"75\u00d745=3375\"],\n  [\"61\u00d795=5795\", \"58\u00d747=2726\"],\n  [\"61\u00d744=2684\", \"81\u00d768=5508\"],\n  [\"23\u00d783=1909\", \"34\u00d761=2074\"],\n  [\"19\u00d727=513\", \"42\u00d797=4074\"],\n  [\"26\u00d715=390\", \"62\u00d744=2728\"],\n  [\"75\u00d713=975\", \"37\u00d712=444\"],\n  [\"77\u00d718=1386\", \"89\u00d785=7565\"],\n  [\"36\u00d760=2160\", \"51\u00d797=4947\"],\n  [\"14\u00d717=238\", \"58\u00d720=1160\"],\n  [\"11\u00d740=440\", \"65\u00d729=1885\"],\n  [\"89\u00d713=1157\", \"47\u00d732=1504\"],\n  [\"83\u00d747=3901\", \"94\u00d777=7238\"],\n  [\"89\u00d754=4806\", \"15\u00d754=810\"],\n  [\"33\u00d787=2871\", \"77\u00d776=5852\"],\n  [\"71\u00d775=5325\", \"22\u00d7100=2200\"],\n  [\"56\u00d737=2072\", \"97\u00d755=5335\"],\n  [\"46\u00d774=3404\", \"89\u00d759=5251\"],\n  [\"98\u00d752=5096\", \"54\u00d755=2970\"],\n  [\"74\u00d787=6438\", \"11\u00d735=385\"],\n  [\"26\u00d726=676\", \"50\u00d771=3550\"],\n  [\"42\u00d760=2520\", \"30\u00d777=2310\"],\n  [\"11\u00d787=957\", \"18\u00d710=180\"],\n  [\"24\u00d727=648\", \"31\u00d726=806\"],\n  [\"79\u00d741=3239\", \"37\u00d759=2183\"],\n  [\"85\u00d799=8415\", \"61\u00d721=1281\"],\n  [\"25\u00d788=2200\", \"52\u00d754=2808\"],\n  [\"76\u00d785=6460\", \"38\u00d737=1406\"],\n  [\"43\u00d763=2709\", \"31\u00d738=1178\"],\n  [\"34\u00d755=1870\", \"61\u00d722=1342\"],\n  [\"54\u00d726=1404\", \"32\u00d736=1152\"],\n  [\"81\u00d736=2916\", \"16\u00d723=368\"],\n  [\"78\u00d715=1170\", \"29\u00d766=1914\"],\n  [\"51\u00d761=3111\", \"96\u00d710=960\"],\n  [\"33\u00d710=330\", \"74\u00d7100=7400\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the 100 multiplication-answer cells\n# with their updated values, in document order. Each old value is unique\n# in the document, so a simple Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n$pairs = @(\n    @(\"2023-06-20 Tuesday\", \"2023-06-21 Wednesday\"),\n    @(\"61\u00d747=2867\", \"17\u00d743=731\"),\n    @(\"65\u00d752=3380\", \"83\u00d786=7138\"),\n    @(\"25\u00d746=1150\", \"43\u00d786=3698\"),\n    @(\"63\u00d746=2898\", \"13\u00d792=1196\"),\n    @(\"14\u00d793=1302\", \"65\u00d781=5265\"),\n    @(\"32\u00d743=1376\", \"87\u00d733=2871\"),\n    @(\"96\u00d712=1152\", \"40\u00d714=560\"),\n    @(\"14\u00d782=1148\", \"33\u00d728=924\"),\n    @(\"86\u00d727=2322\", \"17\u00d766=1122\"),\n    @(\"50\u00d737=1850\", \"37\u00d784=3108\"),\n    @(\"10\u00d779=790\", \"96\u00d798=9408\"),\n    @(\"44\u00d727=1188\", \"97\u00d748=4656\"),\n    @(\"71\u00d725=1775\", \"59\u00d720=1180\"),\n    @(\"12\u00d719=228\", \"46\u00d741=1886\"),\n    @(\"70\u00d790=6300\", \"68\u00d759=4012\"),\n    @(\"32\u00d765=2080\", \"94\u00d772=6768\"),\n    @(\"76\u00d772=5472\", \"32\u00d710=320\"),\n    @(\"45\u00d758=2610\", \"95\u00d780=7600\"),\n    @(\"23\u00d778=1794\", \"70\u00d787=6090\"),\n    @(\"60\u00d770=4200\", \"78\u00d788=6864\"),\n    @(\"51\u00d774=3774\", \"82\u00d760=4920\"),\n    @(\"43\u00d749=2107\", \"16\u00d794=1504\"),\n    @(\"23\u00d734=782\", \"44\u00d711=484\"),\n    @(\"80\u00d758=4640\", \"100\u00d732=3200\"),\n    @(\"72\u00d767=4824\", \"79\u00d753=4187\"),\n    @(\"18\u00d787=1566\", \"13\u00d758=754\"),\n    @(\"78\u00d769=5382\", \"11\u00d762=682\"),\n    @(\"75\u00d744=3300\", \"73\u00d724=1752\"),\n    @(\"84\u00d721=1764\", \"38\u00d7100=3800\"),\n    @(\"17\u00d759=1003\", \"16\u00d780=1280\"),\n    @(\"25\u00d737=925\", \"37\u00d711=407\"),\n    @(\"22\u00d788=1936\", \"74\u00d728=2072\"),\n    @(\"14\u00d739=546\", \"81\u00d727=2187\"),\n    @(\"79\u00d719=1501\", \"51\u00d731=1581\"),\n    @(\"64\u00d757=3648\", \"34\u00d758=1972\"),\n    @(\"53\u00d798=5194\", \"64\u00d799=6336\"),\n    @(\"78\u00d712=936\", \"65\u00d742=2730\"),\n    @(\"18\u00d788=1584\", \"56\u00d750=2800\"),\n    @(\"97\u00d785=8245\", \"80\u00d731=2480\"),\n    @(\"86\u00d763=5418\", \"42\u00d780=3360\"),\n    @(\"12\u00d711=132\", \"94\u00d777=7238\"),\n    @(\"29\u00d798=2842\", \"16\u00d742=672\"),\n    @(\"40\u00d767=2680\", \"20\u00d723=460\"),\n    @(\"93\u00d775=6975\", \"25\u00d715=375\"),\n    @(\"82\u00d742=3444\", \"42\u00d771=2982\"),\n    @(\"26\u00d756=1456\", \"88\u00d783=7304\"),\n    @(\"73\u00d735=2555\", \"59\u00d778=4602\"),\n    @(\"46\u00d768=3128\", \"37\u00d724=888\"),\n    @(\"94\u00d768=6392\", \"24\u00d733=792\"),\n    @(\"93\u00d716=1488\", \"98\u00d729=2842\"),\n    @(\"64\u00d714=896\", \"93\u00d765=6045\"),\n    @(\"16\u00d747=752\", \"38\u00d758=2204\"),\n    @(\"52\u00d792=4784\", \"89\u00d776=6764\"),\n    @(\"83\u00d733=2739\", \"69\u00d747=3243\"),\n    @(\"33\u00d766=2178\", \"36\u00d797=3492\"),\n    @(\"15\u00d732=480\", \"30\u00d716=480\"),\n    @(\"66\u00d789=5874\", \"39\u00d769=2691\"),\n    @(\"19\u00d743=817\", \"48\u00d799=4752\"),\n    @(\"24\u00d776=1824\", \"85\u00d728=2380\"),\n    @(\"20\u00d777=1540\", \"10\u00d735=350\"),\n    @(\"28\u00d721=588\", \"65\u00d755=3575\"),\n    @(\"27\u00d719=513\", \"56\u00d772=4032\"),\n    @(\"20\u00d788=1760\", \"20\u00d718=360\"),\n    @(\"16\u00d734=544\", \"10\u00d798=980\"),\n    @(\"87\u00d732=2784\", \"58\u00d749=2842\"),\n    @(\"40\u00d749=1960\", \"75\u00d745=3375\"),\n    @(\"61\u00d795=5795\", \"58\u00d747=2726\"),\n    @(\"61\u00d744=2684\", \"81\u00d768=5508\"),\n    @(\"23\u00d783=1909\", \"34\u00d761=2074\"),\n    @(\"19\u00d727=513\", \"42\u00d797=4074\"),\n    @(\"26\u00d715=390\", \"62\u00d744=2728\"),\n    @(\"75\u00d713=975\", \"37\u00d712=444\"),\n    @(\"77\u00d718=1386\", \"89\u00d785=7565\"),\n    @(\"36\u00d760=2160\", \"51\u00d797=4947\"),\n    @(\"14\u00d717=238\", \"58\u00d720=1160\"),\n    @(\"11\u00d740=440\", \"65\u00d729=1885\"),\n    @(\"89\u00d713=1157\", \"47\u00d732=1504\"),\n    @(\"83\u00d747=3901\", \"94\u00d777=7238\"),\n    @(\"89\u00d754=4806\", \"15\u00d754=810\"),\n    @(\"33\u00d787=2871\", \"77\u00d776=5852\"),\n    @(\"71\u00d775=5325\", \"22\u00d7100=2200\"),\n    @(\"56\u00d737=2072\", \"97\u00d755=5335\"),\n    @(\"46\u00d774=3404\", \"89\u00d759=5251\"),\n    @(\"98\u00d752=5096\", \"54\u00d755=2970\"),\n    @(\"74\u00d787=6438\", \"11\u00d735=385\"),\n    @(\"26\u00d726=676\", \"50\u00d771=3550\"),\n    @(\"42\u00d760=2520\", \"30\u00d777=2310\"),\n    @(\"11\u00d787=957\", \"18\u00d710=180\"),\n    @(\"24\u00d727=648\", \"31\u00d726=806\"),\n    @(\"79\u00d741=3239\", \"37\u00d759=2183\"),\n    @(\"85\u00d799=8415\", \"61\u00d721=1281\"),\n    @(\"25\u00d788=2200\", \"52\u00d754=2808\"),\n    @(\"76\u00d785=6460\", \"38\u00d737=1406\"),\n    @(\"43\u00d763=2709\", \"31\u00d738=1178\"),\n    @(\"34\u00d755=1870\", \"61\u00d722=1342\"),\n    @(\"54\u00d726=1404\", \"32\u00d736=1152\"),\n    @(\"81\u00d736=2916\", \"16\u00d723=368\"),\n    @(\"78\u00d715=1170\", \"29\u00d766=1914\"),\n    @(\"51\u00d761=3111\", \"96\u00d710=960\"),\n    @(\"33\u00d710=330\", \"74\u00d7100=7400\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $result) {\n        Write-Output \"FAILED: $old -> $new\"\n    }\n}\n"}
